$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 3 (swap in former row 5 price data, keep date shift) ---
$ws.Range("D3").Value = 44160
$ws.Range("N3").Value = 19000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 19500
$ws.Range("S3").Value = 1500

# --- Update row 4 (becomes former row 3's data, quality/unit change to Segunda) ---
$ws.Range("D4").Value = 44167
$ws.Range("L4").Value = "Segunda"
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 19000
$ws.Range("P4").Value = 18500
$ws.Range("Q4").Value = "$/caja 13 kilos"
$ws.Range("S4").Value = 1423
$ws.Range("T4").Value = 13

# --- Update row 5 (new "Especial" quality entry) ---
$ws.Range("D5").Value = 44475
$ws.Range("L5").Value = "Especial"
$ws.Range("N5").Value = 32000
$ws.Range("O5").Value = 33000
$ws.Range("P5").Value = 32500
$ws.Range("Q5").Value = "$/caja 12 kilos"
$ws.Range("S5").Value = 2708
$ws.Range("T5").Value = 12

# --- Add new row 6, a copy of the original row 4 data (Primera / $/bandeja 10 kilos) ---
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C6").Value = "Arica y Parinacota"
$ws.Range("D6").Value = 44468
$ws.Range("D6").NumberFormat = $ws.Range("D5").NumberFormat
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100107
$ws.Range("H6").Value = "Otros"
$ws.Range("I6").Value = 100107002
$ws.Range("J6").Value = "Chirimoya"
$ws.Range("K6").Value = "Cultivar IV Región"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 29000
$ws.Range("O6").Value = 30000
$ws.Range("P6").Value = 29500
$ws.Range("Q6").Value = "$/bandeja 10 kilos"
$ws.Range("R6").Value = "Región de Coquimbo"
$ws.Range("S6").Value = 2950
$ws.Range("T6").Value = 10
